# Generate Report for Handback
# Updates the localization-status workbook to reflect that both files have
# been handed back and are back in sync with en-US: the Status column
# moves from "Ready for handoff" to "Handed back: in sync with en-US", the
# "Latest Handback DateTime" gets a real timestamp, and the new
# "Latest Target File" / "Latest Handback File" columns (F/G) are filled in
# with hyperlinks mirroring the source/handoff file links.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet - its zh-cn/de-de status cells share the same underlying
# "Ready for handoff" string as the per-locale sheets, so they flip to the
# new status too.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusText
$wsOverview.Range("C2").Value = $statusText
$wsOverview.Range("B3").Value = $statusText
$wsOverview.Range("C3").Value = $statusText

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Status column (C) -> handed back
$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

# Latest Handback DateTime (H)
$wsZh.Range("H2").Value = "2016-03-23 16:16:19"
$wsZh.Range("H3").Value = "2016-03-23 16:16:19"

# New columns: F = Latest Target File, G = Latest Handback File
$wsZh.Range("F2").Value = "b5886bb2-4392-4fb5-aa09-bc3afd145ec7.md"
$wsZh.Range("G2").Value = "b5886bb2-4392-4fb5-aa09-bc3afd145ec7.858c7bc584b9bdbdeecaa789bb0793655fb31591.zh-cn.xlf"
$wsZh.Range("F3").Value = "c8eff952-8675-467a-bd71-096699da26e5.md"
$wsZh.Range("G3").Value = "c8eff952-8675-467a-bd71-096699da26e5.6d1af518bfebeacd791e4e9404f9ec8cd44bda57.zh-cn.xlf"

# Rebuild hyperlinks in display order so relationship ids line up the way
# Excel would naturally emit them (A2, D2, F2, G2, A3, D3, F3, G3).
$wsZh.Hyperlinks.Delete()

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/7914b7dae3fa6a6152685112a99b2f7d7cbab015/e2e/b5886bb2-4392-4fb5-aa09-bc3afd145ec7.md", "", "", "b5886bb2-4392-4fb5-aa09-bc3afd145ec7.md")
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a8ca57f81987fd8ee0baa1798cbfa3b591582658/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/b5886bb2-4392-4fb5-aa09-bc3afd145ec7.858c7bc584b9bdbdeecaa789bb0793655fb31591.zh-cn.xlf", "", "", "b5886bb2-4392-4fb5-aa09-bc3afd145ec7.858c7bc584b9bdbdeecaa789bb0793655fb31591.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/7914b7dae3fa6a6152685112a99b2f7d7cbab015/e2e/b5886bb2-4392-4fb5-aa09-bc3afd145ec7.md", "", "", "b5886bb2-4392-4fb5-aa09-bc3afd145ec7.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a8ca57f81987fd8ee0baa1798cbfa3b591582658/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/b5886bb2-4392-4fb5-aa09-bc3afd145ec7.858c7bc584b9bdbdeecaa789bb0793655fb31591.zh-cn.xlf", "", "", "b5886bb2-4392-4fb5-aa09-bc3afd145ec7.858c7bc584b9bdbdeecaa789bb0793655fb31591.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/7914b7dae3fa6a6152685112a99b2f7d7cbab015/e2e/c8eff952-8675-467a-bd71-096699da26e5.md", "", "", "c8eff952-8675-467a-bd71-096699da26e5.md")
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a8ca57f81987fd8ee0baa1798cbfa3b591582658/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/c8eff952-8675-467a-bd71-096699da26e5.6d1af518bfebeacd791e4e9404f9ec8cd44bda57.zh-cn.xlf", "", "", "c8eff952-8675-467a-bd71-096699da26e5.6d1af518bfebeacd791e4e9404f9ec8cd44bda57.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/7914b7dae3fa6a6152685112a99b2f7d7cbab015/e2e/c8eff952-8675-467a-bd71-096699da26e5.md", "", "", "c8eff952-8675-467a-bd71-096699da26e5.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a8ca57f81987fd8ee0baa1798cbfa3b591582658/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/c8eff952-8675-467a-bd71-096699da26e5.6d1af518bfebeacd791e4e9404f9ec8cd44bda57.zh-cn.xlf", "", "", "c8eff952-8675-467a-bd71-096699da26e5.6d1af518bfebeacd791e4e9404f9ec8cd44bda57.zh-cn.xlf")

# Restore the workbook's original hyperlink look (underlined, #6495ED) on
# every linked cell -- re-adding hyperlinks resets them to Excel's default
# themed hyperlink style, so force them back to match A2/D2/A3/D3's
# pre-existing formatting.
$wsZh.Range("A2").Font.Underline = 2
$wsZh.Range("A2").Font.Color = 15570276
$wsZh.Range("D2").Font.Underline = 2
$wsZh.Range("D2").Font.Color = 15570276
$wsZh.Range("F2").Font.Underline = 2
$wsZh.Range("F2").Font.Color = 15570276
$wsZh.Range("G2").Font.Underline = 2
$wsZh.Range("G2").Font.Color = 15570276
$wsZh.Range("A3").Font.Underline = 2
$wsZh.Range("A3").Font.Color = 15570276
$wsZh.Range("D3").Font.Underline = 2
$wsZh.Range("D3").Font.Color = 15570276
$wsZh.Range("F3").Font.Underline = 2
$wsZh.Range("F3").Font.Color = 15570276
$wsZh.Range("G3").Font.Underline = 2
$wsZh.Range("G3").Font.Color = 15570276

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

$wsDe.Range("H2").Value = "2016-03-23 16:16:34"
$wsDe.Range("H3").Value = "2016-03-23 16:16:34"

$wsDe.Range("F2").Value = "b5886bb2-4392-4fb5-aa09-bc3afd145ec7.md"
$wsDe.Range("G2").Value = "b5886bb2-4392-4fb5-aa09-bc3afd145ec7.858c7bc584b9bdbdeecaa789bb0793655fb31591.de-de.xlf"
$wsDe.Range("F3").Value = "c8eff952-8675-467a-bd71-096699da26e5.md"
$wsDe.Range("G3").Value = "c8eff952-8675-467a-bd71-096699da26e5.6d1af518bfebeacd791e4e9404f9ec8cd44bda57.de-de.xlf"

$wsDe.Hyperlinks.Delete()

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/7914b7dae3fa6a6152685112a99b2f7d7cbab015/e2e/b5886bb2-4392-4fb5-aa09-bc3afd145ec7.md", "", "", "b5886bb2-4392-4fb5-aa09-bc3afd145ec7.md")
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6e6c1d63a22460678607c4d735d55c5b680947fa/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/b5886bb2-4392-4fb5-aa09-bc3afd145ec7.858c7bc584b9bdbdeecaa789bb0793655fb31591.de-de.xlf", "", "", "b5886bb2-4392-4fb5-aa09-bc3afd145ec7.858c7bc584b9bdbdeecaa789bb0793655fb31591.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/7914b7dae3fa6a6152685112a99b2f7d7cbab015/e2e/b5886bb2-4392-4fb5-aa09-bc3afd145ec7.md", "", "", "b5886bb2-4392-4fb5-aa09-bc3afd145ec7.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6e6c1d63a22460678607c4d735d55c5b680947fa/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/b5886bb2-4392-4fb5-aa09-bc3afd145ec7.858c7bc584b9bdbdeecaa789bb0793655fb31591.de-de.xlf", "", "", "b5886bb2-4392-4fb5-aa09-bc3afd145ec7.858c7bc584b9bdbdeecaa789bb0793655fb31591.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/7914b7dae3fa6a6152685112a99b2f7d7cbab015/e2e/c8eff952-8675-467a-bd71-096699da26e5.md", "", "", "c8eff952-8675-467a-bd71-096699da26e5.md")
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6e6c1d63a22460678607c4d735d55c5b680947fa/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/c8eff952-8675-467a-bd71-096699da26e5.6d1af518bfebeacd791e4e9404f9ec8cd44bda57.de-de.xlf", "", "", "c8eff952-8675-467a-bd71-096699da26e5.6d1af518bfebeacd791e4e9404f9ec8cd44bda57.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/7914b7dae3fa6a6152685112a99b2f7d7cbab015/e2e/c8eff952-8675-467a-bd71-096699da26e5.md", "", "", "c8eff952-8675-467a-bd71-096699da26e5.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6e6c1d63a22460678607c4d735d55c5b680947fa/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/c8eff952-8675-467a-bd71-096699da26e5.6d1af518bfebeacd791e4e9404f9ec8cd44bda57.de-de.xlf", "", "", "c8eff952-8675-467a-bd71-096699da26e5.6d1af518bfebeacd791e4e9404f9ec8cd44bda57.de-de.xlf")

# Restore the workbook's original hyperlink look (underlined, #6495ED) on
# every linked cell -- re-adding hyperlinks resets them to Excel's default
# themed hyperlink style, so force them back to match A2/D2/A3/D3's
# pre-existing formatting.
$wsDe.Range("A2").Font.Underline = 2
$wsDe.Range("A2").Font.Color = 15570276
$wsDe.Range("D2").Font.Underline = 2
$wsDe.Range("D2").Font.Color = 15570276
$wsDe.Range("F2").Font.Underline = 2
$wsDe.Range("F2").Font.Color = 15570276
$wsDe.Range("G2").Font.Underline = 2
$wsDe.Range("G2").Font.Color = 15570276
$wsDe.Range("A3").Font.Underline = 2
$wsDe.Range("A3").Font.Color = 15570276
$wsDe.Range("D3").Font.Underline = 2
$wsDe.Range("D3").Font.Color = 15570276
$wsDe.Range("F3").Font.Underline = 2
$wsDe.Range("F3").Font.Color = 15570276
$wsDe.Range("G3").Font.Underline = 2
$wsDe.Range("G3").Font.Color = 15570276
